# "Generate Report for handback"
# Refresh the handoff/handback timestamps recorded for the two
# xlf files ("...zh-cn.xlf" and "...de-de.xlf") on the zh-cn / de-de
# report sheets. Only the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) on row 2 of each language
# sheet are refreshed with newly generated timestamps - everything else
# (row 3, other columns, the Overview sheet) is left untouched.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-13 02:03:43"
$wsZhCn.Range("G2").Value = "2016-01-13 02:05:00"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-13 02:04:07"
$wsDeDe.Range("G2").Value = "2016-01-13 02:05:32"
